$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 733
$ws.Range("I2").Value = 599.5
$ws.Range("K2").Value = 599.5
$ws.Range("M2").Value = -486.5
$ws.Range("H17").Value = 3383.3333
$ws.Range("H33").Value = 582.0769
$ws.Range("I33").Value = 582.0769
$ws.Range("K33").Value = 582.0769
$ws.Range("M33").Value = -353.0769
$ws.Range("H40").Value = 2299.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2299.5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2299.5
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2649.5
$ws.Range("H61").Value = 608.75
$ws.Range("I61").Value = 511.66666
$ws.Range("K61").Value = 1534.99998
$ws.Range("M61").Value = -1362.99998
$ws.Range("H98").Value = 1132.3334
$ws.Range("I98").Value = 1059.2
$ws.Range("K98").Value = 1059.2
$ws.Range("M98").Value = 438.8
$ws.Range("H122").Value = 1132.3334
$ws.Range("I122").Value = 1059.2
$ws.Range("K122").Value = 3177.6
$ws.Range("M122").Value = -727.6000000000004
$ws.Range("H123").Value = 130000
$ws.Range("J123").Value = 130000
$ws.Range("L123").Value = 130000
$ws.Range("N123").Value = -139800
$ws.Range("H126").Value = 99999
$ws.Range("J126").Value = 99999
$ws.Range("L126").Value = 99999
$ws.Range("N126").Value = -109879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("H105").Value = 1874.8572
$ws.Range("I105").Value = 1874.8572
$ws.Range("K105").Value = 1874.8572
$ws.Range("M105").Value = -127.8571999999999
$ws.Range("H107").Value = 3636.4285
$ws.Range("I107").Value = 3553.2
$ws.Range("K107").Value = 3553.2
$ws.Range("M107").Value = -1633.2
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 33333.332
$ws.Range("I44").Value = 31500
$ws.Range("K44").Value = 31500
$ws.Range("M44").Value = -31058
$ws.Range("H99").Value = 4666.6665
$ws.Range("I99").Value = 4666.6665
$ws.Range("K99").Value = 4666.6665
$ws.Range("M99").Value = -3168.6665
$ws.Range("H105").Value = 587
$ws.Range("I105").Value = 498.16666
$ws.Range("J105").Value = 764.6667
$ws.Range("K105").Value = 498.16666
$ws.Range("L105").Value = 764.6667
$ws.Range("M105").Value = 1248.83334
$ws.Range("N105").Value = -4258.6667
$ws.Range("H126").Value = 4666.6665
$ws.Range("I126").Value = 4666.6665
$ws.Range("K126").Value = 13999.9995
$ws.Range("M126").Value = -11529.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6966.3335
$ws.Range("I5").Value = 6966.3335
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 20899.0005
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -20787.0005
$ws.Range("N5").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H117").Value = 2612.2
$ws.Range("I117").Value = 765.25
$ws.Range("K117").Value = 2295.75
$ws.Range("M117").Value = 1146.25
$ws.Range("H131").Value = 2299.7144
$ws.Range("J131").Value = 2299.7144
$ws.Range("L131").Value = 6899.1432
$ws.Range("N131").Value = -16979.1432
$ws.Range("H132").Value = 617.4286
$ws.Range("J132").Value = 924.5
$ws.Range("L132").Value = 8320.5
$ws.Range("N132").Value = -13380.5
$ws.Range("H135").Value = 6966.3335
$ws.Range("I135").Value = 6966.3335
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 62697.0015
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -60162.0015
$ws.Range("N135").ClearContents()
$ws.Range("H140").Value = 2349
$ws.Range("I140").Value = 2349
$ws.Range("K140").Value = 7047
$ws.Range("M140").Value = -1867

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 21000
$ws.Range("J40").Value = 21000
$ws.Range("L40").Value = 21000
$ws.Range("N40").Value = -21302
$ws.Range("H102").Value = 1508.5
$ws.Range("I102").Value = 1542.7778
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 1542.7778
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = 79.22219999999993
$ws.Range("N102").Value = -4444
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8887.241
$ws.Range("I7").Value = 8619.629999999999
$ws.Range("J7").Value = 12500
$ws.Range("K7").Value = 8619.629999999999
$ws.Range("L7").Value = 12500
$ws.Range("M7").Value = -8507.629999999999
$ws.Range("N7").Value = -12724
$ws.Range("H54").Value = 31542
$ws.Range("J54").Value = 31542
$ws.Range("L54").Value = 31542
$ws.Range("N54").Value = -32830
$ws.Range("H61").Value = 3099.875
$ws.Range("I61").Value = 2971.2856
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2971.2856
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2769.2856
$ws.Range("N61").Value = -4404
$ws.Range("H113").Value = 3099.875
$ws.Range("I113").Value = 2971.2856
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2971.2856
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -801.2856000000002
$ws.Range("N113").Value = -8340
$ws.Range("H126").Value = 8887.241
$ws.Range("I126").Value = 8619.629999999999
$ws.Range("J126").Value = 12500
$ws.Range("K126").Value = 25858.89
$ws.Range("L126").Value = 37500
$ws.Range("M126").Value = -23388.89
$ws.Range("N126").Value = -42440
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H122").Value = 1391.5714
$ws.Range("I122").Value = 936
$ws.Range("K122").Value = 2808
$ws.Range("M122").Value = -358
$ws.Range("H126").Value = 4425.6924
$ws.Range("I126").Value = 3518
$ws.Range("K126").Value = 10554
$ws.Range("M126").Value = -8084
$ws.Range("H132").Value = 3998.5
$ws.Range("I132").Value = 3998.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11995.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9465.5
$ws.Range("N132").ClearContents()
